$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank columns before K so the existing "Optie onbeperkt reizen" (K:L)
# table shifts to N:O and the existing "Auto kopen" (N:O) table shifts to Q:R,
# exactly like Excel does when you insert columns (formulas auto-adjust).
$ws.Range("K1:M1").EntireColumn.Insert()

# Re-assert these two literals (the engine reserialises them with float noise
# across the column-insert shift; they are plain literals, not formulas).
$ws.Range("R4").Value = 7.65
$ws.Range("R5").Value = 1.65

# --- New "Optie boete" table in the freed-up K:L columns ---
# (same row labels as the "Optie geen kaart" table in B:C, new header + values)
$ws.Range("K2").Value = "Optie boete"

$ws.Range("K3").Value = "Kaart"
$ws.Range("L3").Formula = "=2*97"

$ws.Range("K4").Value = "Korting student"
$ws.Range("L4").Value = 0

$ws.Range("K7").Value = "# Ritten voor 9 uur"
$ws.Range("L7").Value = 20

$ws.Range("K8").Value = "korting spits"
$ws.Range("L8").Value = 1

$ws.Range("K9").Value = "prijs NS"
$ws.Range("L9").Formula = "=C9*(1-L8)"

$ws.Range("K10").Value = "prijs Arriva"
$ws.Range("L10").Formula = "=C10*(1-L8)"

$ws.Range("K11").Value = "Ritten buiten spits"
$ws.Range("L11").Value = 6

$ws.Range("K12").Value = "korting dal"
$ws.Range("L12").Value = 1

$ws.Range("K13").Value = "prijs NS"
$ws.Range("L13").Formula = "=C13*(1-L12)"

$ws.Range("K14").Value = "prijs Arriva"
$ws.Range("L14").Formula = "=C14*(1-L12)"

$ws.Range("K16").Value = "Kosten wachttijd"
$ws.Range("L16").Formula = "=(L7+L11)*0.5*10"

$ws.Range("K17").Value = "totaal OV kosten Damen"
$ws.Range("L17").Formula = "=L7*(L9+L10)+L11*(L13+L14)+L3+L4"

$ws.Range("K18").Value = "overige OV kosten"
$ws.Range("L18").Value = 30

$ws.Range("K20").Value = "totaal"
$ws.Range("L20").Formula = "=L18+L17+L16"

# --- "Goedkoopst" marker, formatted with the built-in "Good" cell style ---
$ws.Range("K22").Value = "Goedkoopst"
$ws.Range("K22:L22").Style = "Good"

$ws.Range("K23").Select()
